$d = $word.ActiveDocument

# Remove the sentence "The full model could potentially predict up to a
# 2.26-point differential in negative affect." (and the trailing space
# that followed it), which previously sat between "... in negative
# affect. " and "Thus, as hypothesized, ...". Replacing the exact
# (unique) run of text with an empty string deletes those four runs
# outright while leaving the surrounding text untouched.
$rng = $d.Content
$found = $rng.Find.Execute(
    "The full model could potentially predict up to a 2.26-point differential in negative affect. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

if (-not $found) {
    throw "Target sentence not found - document may already be edited or text changed."
}
